# edit.ps1 - applies the WBA formatting/content changes described in the
# commit: column width tweaks, "Testing & Reviewing" grammar-check run
# split, "More Allies and Enemies" / "Shyam and Oskar" re-wrapped onto two
# lines each, "item block" -> "mystery block" (split into separate runs),
# and removal of the trailing blank paragraph at the end of the document.

$d = $word.ActiveDocument

function Set-ParagraphXml($rng, [string]$paragraphsXml) {
    # Replaces the OOXML of the paragraph(s) touched by $rng with the raw
    # <w:p>...</w:p> markup supplied in $paragraphsXml. InsertXML operates
    # at paragraph granularity, so $paragraphsXml must be a self-contained
    # replacement for the whole paragraph(s) that $rng currently overlaps.
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $paragraphsXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1) Column widths: first column 2689 -> 2830 dxa, second column
#    1991 -> 1850 dxa (dxa / 20 = points). Setting the width of one cell
#    re-flows the whole column (gridCol + every tcW in that column).
# ---------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Width = 2830 / 20
$t.Cell(1, 2).Width = 1850 / 20

# ---------------------------------------------------------------------
# 2) "Testing & Reviewing" header cell: split into two runs with
#    proofErr gramStart/gramEnd markers bracketing "Reviewing".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Testing & Reviewing") | Out-Null
$xml = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Testing &amp; </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Reviewing</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
Set-ParagraphXml $rng $xml

# ---------------------------------------------------------------------
# 3) "2: More Allies and Enemies" cell: keep the "2: " run, but put
#    "More " and "Allies and Enemies" into two separate paragraphs.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("More Allies and Enemies") | Out-Null
$xml = '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">2: </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">More </w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Allies and Enemies</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $rng $xml

# ---------------------------------------------------------------------
# 4) "Shyam and Oskar" cell: "Shyam " on its own paragraph (kept as two
#    runs: "Shyam" + " "), "and Oskar" on the next paragraph.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Shyam and Oskar") | Out-Null
$xml = '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Shyam</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>and Oskar</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $rng $xml

# ---------------------------------------------------------------------
# 5) "5 C2: - Random trees, enemies, item block" cell: keep "5 ", "C2: "
#    and "-" runs, change "item" to "mystery" and split the tail into
#    three runs: ", enemies,", " mystery" and " block".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Random trees, enemies, item block") | Out-Null
$xml = '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">5 </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">C2: </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>-</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Random trees, enemies,</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> mystery</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> block</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $rng $xml

# ---------------------------------------------------------------------
# 6) Remove the trailing empty paragraph after "I agree to this WBA -
#    Vedesh" (just before the sectPr).
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($n - 1)
$last = $d.Paragraphs.Item($n)
$d.Range($secondLast.Range.End - 1, $last.Range.End).Delete()
